$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = $null

$ws.Range("F14").Value = 11.6
$ws.Range("J14").Value = 11.6
$ws.Range("F15").Value = 11.6
$ws.Range("F16").Value = 29
$ws.Range("F17").Value = 52.6
$ws.Range("F18").Value = 64.5

$ws.Range("C3:G3").Select()
